$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the formatting of the last existing data row (145) down into the
# three new rows (146:148) so the date column keeps the same short-date
# number format (style index) as the rest of the table.
$ws.Range("A145:M145").Copy()
$ws.Range("A146:M148").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 146: Eric / Workout / Agile Antelope
$ws.Cells.Item(146, 1).Value = "Eric"
$ws.Cells.Item(146, 2).Value = 45475
$ws.Cells.Item(146, 3).Value = "Workout"
$ws.Cells.Item(146, 4).Value = 76
$ws.Cells.Item(146, 5).Value = 0
$ws.Cells.Item(146, 6).Value = 0
$ws.Cells.Item(146, 7).Value = 21
$ws.Cells.Item(146, 8).Value = 46
$ws.Cells.Item(146, 9).Value = 9
$ws.Cells.Item(146, 10).Value = 1
$ws.Cells.Item(146, 11).Value = 0
$ws.Cells.Item(146, 12).Value = "Agile Antelope"
$ws.Cells.Item(146, 13).Value = 4

# Row 147: Eric / Walk / Agile Antelope
$ws.Cells.Item(147, 1).Value = "Eric"
$ws.Cells.Item(147, 2).Value = 45475
$ws.Cells.Item(147, 3).Value = "Walk"
$ws.Cells.Item(147, 4).Value = 23
$ws.Cells.Item(147, 5).Value = 1.23
$ws.Cells.Item(147, 6).Value = 46
$ws.Cells.Item(147, 7).Value = 23
$ws.Cells.Item(147, 8).Value = 0
$ws.Cells.Item(147, 9).Value = 0
$ws.Cells.Item(147, 10).Value = 0
$ws.Cells.Item(147, 11).Value = 0
$ws.Cells.Item(147, 12).Value = "Agile Antelope"
$ws.Cells.Item(147, 13).Value = 4

# Row 148: Jeremiah / Workout / Agile Antelope
$ws.Cells.Item(148, 1).Value = "Jeremiah"
$ws.Cells.Item(148, 2).Value = 45475
$ws.Cells.Item(148, 3).Value = "Workout"
$ws.Cells.Item(148, 4).Value = 79
$ws.Cells.Item(148, 5).Value = 0
$ws.Cells.Item(148, 6).Value = 0
$ws.Cells.Item(148, 7).Value = 15
$ws.Cells.Item(148, 8).Value = 38
$ws.Cells.Item(148, 9).Value = 25
$ws.Cells.Item(148, 10).Value = 2
$ws.Cells.Item(148, 11).Value = 0
$ws.Cells.Item(148, 12).Value = "Agile Antelope"
$ws.Cells.Item(148, 13).Value = 4

# Match the saved view state: the active cell moves to the row right after
# the newly appended data, same as the source workbook.
$ws.Range("A149").Select()
